$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

# --- Crime statistics table updates (rows 14-33) ---
# Row 14
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = -50
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -92.857142857142

# Row 15
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -7.692307692307
$ws.Range("N15").Value = 71.428571428571

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 31
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = -16.216216216216
$ws.Range("I16").Value = 138
$ws.Range("J16").Value = 181
$ws.Range("K16").Value = -23.756906077348
$ws.Range("L16").Value = -2.81690140845
$ws.Range("M16").Value = -4.827586206896
$ws.Range("N16").Value = -77.227722772277

# Row 17
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 5.555555555555
$ws.Range("F17").Value = 58
$ws.Range("G17").Value = 59
$ws.Range("H17").Value = -1.694915254237
$ws.Range("I17").Value = 234
$ws.Range("J17").Value = 283
$ws.Range("K17").Value = -17.314487632508
$ws.Range("L17").Value = -2.904564315352
$ws.Range("M17").Value = 120.754716981132
$ws.Range("N17").Value = 47.169811320754

# Row 18
$ws.Range("C18").Value = "0"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 27.272727272727
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 109
$ws.Range("K18").Value = -35.779816513761
$ws.Range("L18").Value = -16.666666666666
$ws.Range("M18").Value = -50.354609929078
$ws.Range("N18").Value = -92.099322799097

# Row 19
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 5.555555555555
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -20.27027027027
$ws.Range("I19").Value = 269
$ws.Range("J19").Value = 380
$ws.Range("K19").Value = -29.210526315789
$ws.Range("L19").Value = -20.648967551622
$ws.Range("M19").Value = 35.858585858585
$ws.Range("N19").Value = -35.024154589372

# Row 20
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 11.111111111111
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -14.814814814814
$ws.Range("I20").Value = 80
$ws.Range("J20").Value = 101
$ws.Range("K20").Value = -20.79207920792
$ws.Range("L20").Value = -24.528301886792
$ws.Range("M20").Value = 31.147540983606
$ws.Range("N20").Value = -89.610389610389

# Row 21
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -6.896551724137
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 210
$ws.Range("H21").Value = -11.428571428571
$ws.Range("I21").Value = 804
$ws.Range("J21").Value = 1069
$ws.Range("K21").Value = -24.789522918615
$ws.Range("L21").Value = -12.987012987013
$ws.Range("M21").Value = 21.450151057401
$ws.Range("N21").Value = -71.848739495798

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 19
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 11.764705882352
$ws.Range("L22").Value = 18.75
$ws.Range("M22").Value = 46.153846153846

# Row 24
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 79
$ws.Range("E24").Value = -41.772151898734
$ws.Range("F24").Value = 183
$ws.Range("G24").Value = 254
$ws.Range("H24").Value = -27.952755905511
$ws.Range("I24").Value = 919
$ws.Range("J24").Value = 1249
$ws.Range("K24").Value = -26.421136909527
$ws.Range("L24").Value = -11.378977820636
$ws.Range("M24").Value = 60.664335664335

# Row 25
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 60
$ws.Range("E25").Value = -58.333333333333
$ws.Range("F25").Value = 113
$ws.Range("G25").Value = 201
$ws.Range("H25").Value = -43.781094527363
$ws.Range("I25").Value = 690
$ws.Range("J25").Value = 1005
$ws.Range("K25").Value = -31.343283582089
$ws.Range("L25").Value = -10.737386804657

# Row 26
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 102
$ws.Range("G26").Value = 103
$ws.Range("H26").Value = -0.970873786407
$ws.Range("I26").Value = 395
$ws.Range("J26").Value = 467
$ws.Range("K26").Value = -15.417558886509
$ws.Range("L26").Value = -2.227722772277
$ws.Range("M26").Value = 74.008810572687

# Row 27
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 6.666666666666

# Row 28
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 9
$ws.Range("H28").Value = 125
$ws.Range("I28").Value = 43
$ws.Range("K28").Value = 10.25641025641
$ws.Range("L28").Value = -14

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("F29").Value = 4
$ws.Range("I29").Value = 5
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = 150
$ws.Range("N29").Value = -80.76923076923

# Row 30
$ws.Range("C30").Value = 2
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 3
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = -88

# Row 33
$ws.Range("G33").Value = "0"
$ws.Range("H33").Value = "***.*"
